$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-11-23 Saturday" "2024-11-24 Sunday"
Replace-Text "218÷7=" "152÷2="
Replace-Text "540÷3=" "303÷2="
Replace-Text "165÷4=" "485÷8="
Replace-Text "249÷9=" "752÷8="
Replace-Text "804÷9=" "429÷9="
Replace-Text "985÷8=" "170÷7="
Replace-Text "344÷9=" "388÷7="
Replace-Text "221÷2=" "291÷2="
Replace-Text "134÷2=" "869÷4="
Replace-Text "192÷9=" "494÷2="
Replace-Text "456÷3=" "586÷6="
Replace-Text "910÷5=" "776÷6="
Replace-Text "543÷8=" "612÷3="
Replace-Text "286÷5=" "711÷9="
Replace-Text "576÷7=" "318÷4="
Replace-Text "655÷4=" "155÷6="
Replace-Text "480÷7=" "135÷8="
Replace-Text "939÷7=" "848÷8="
Replace-Text "863÷7=" "762÷8="
Replace-Text "270÷2=" "145÷2="
Replace-Text "172÷9=" "281÷7="
Replace-Text "913÷5=" "757÷5="
Replace-Text "909÷7=" "512÷5="
Replace-Text "675÷7=" "995÷5="
Replace-Text "128÷9=" "407÷9="
